$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy cell formatting from the (now-shifted) original D:E columns (now F:G) into the new D:E columns.
# Done per contiguous data block so the fully-blank separator rows (36, 78) and the
# section-header-only rows (37, 79) are left untouched, matching the source layout.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

# Populate the two new columns with the latest two quarters of data
$ws.Range("D7").Value = 43463
$ws.Range("E7").Value = 43372
$ws.Range("D8").Value = 77000
$ws.Range("E8").Value = 76600
$ws.Range("D9").Value = 35000
$ws.Range("E9").Value = 33000
$ws.Range("D10").Value = 42000
$ws.Range("E10").Value = 43600
$ws.Range("D12").Value = 12800
$ws.Range("E12").Value = 12700
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "NA"
$ws.Range("D17").Value = 66100
$ws.Range("E17").Value = 62300
$ws.Range("D18").Value = 10900
$ws.Range("E18").Value = 14300
$ws.Range("D20").Value = 900
$ws.Range("E20").Value = 300
$ws.Range("D21").Value = 13600
$ws.Range("E21").Value = 16200
$ws.Range("D22").Value = 100
$ws.Range("E22").Value = 100
$ws.Range("D23").Value = 11600
$ws.Range("E23").Value = 14500
$ws.Range("D24").Value = -400
$ws.Range("E24").Value = 2900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 12000
$ws.Range("E26").Value = 11600
$ws.Range("D27").Value = 12000
$ws.Range("E27").Value = 11600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -900
$ws.Range("E32").Value = -300
$ws.Range("D33").Value = 12000
$ws.Range("E33").Value = 11600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 12000
$ws.Range("E35").Value = 11600
$ws.Range("D38").Value = 43463
$ws.Range("E38").Value = 43372
$ws.Range("D41").Value = 111000
$ws.Range("E41").Value = 96100
$ws.Range("D42").Value = 40800
$ws.Range("E42").Value = 76800
$ws.Range("D43").Value = 50900
$ws.Range("E43").Value = 50300
$ws.Range("D44").Value = 62100
$ws.Range("E44").Value = 60700
$ws.Range("D45").Value = 6100
$ws.Range("E45").Value = 7100
$ws.Range("D46").Value = 270900
$ws.Range("E46").Value = 291100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 47900
$ws.Range("E48").Value = 43100
$ws.Range("D49").Value = 53700
$ws.Range("E49").Value = 14700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 3200
$ws.Range("E52").Value = 7500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 375600
$ws.Range("E54").Value = 356500
$ws.Range("D57").Value = 16500
$ws.Range("E57").Value = 20900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 43200
$ws.Range("E59").Value = 38200
$ws.Range("D60").Value = 59800
$ws.Range("E60").Value = 59100
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 3000
$ws.Range("E62").Value = 3300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 62800
$ws.Range("E66").Value = 62400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 67400
$ws.Range("E72").Value = 55400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 312900
$ws.Range("E76").Value = 294100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("E80").Value = 43372
$ws.Range("D81").Value = 12000
$ws.Range("E81").Value = 11600
$ws.Range("D83").Value = 1900
$ws.Range("E83").Value = 1600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 19900
$ws.Range("E89").Value = 24400
$ws.Range("D91").Value = -4500
$ws.Range("E91").Value = -1200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -4700
$ws.Range("E94").Value = -15300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -500
$ws.Range("E100").Value = 700
$ws.Range("D101").Value = 100
$ws.Range("E101").Value = 100
$ws.Range("D102").Value = 14800
$ws.Range("E102").Value = 9900
